# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the row corresponding to the
# 2c0b003b-fab4-49a4-97a0-24c304bafa67 file's dependency (0f2df5a2 row) on both the
# zh-cn and de-de localization-status sheets, reflecting the newly generated handoff.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-09 02:56:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-09 02:57:08"
